$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.044402855926531
$ws.Cells.Item(2, 4).Value = 1.03940056114219
$ws.Cells.Item(2, 5).Value = 1.04216980029393
$ws.Cells.Item(2, 6).Value = 1.052252153643817
$ws.Cells.Item(2, 9).Value = 1.042790746720722
$ws.Cells.Item(2, 10).Value = 1.049468130274581
$ws.Cells.Item(2, 11).Value = 1.042185855742719
$ws.Cells.Item(2, 12).Value = 1.044947245901242
$ws.Cells.Item(2, 13).Value = 1.055001400733993
$ws.Cells.Item(2, 14).Value = 1.050958495395248

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.045951135086813
$ws.Cells.Item(3, 4).Value = 1.040130348693367
$ws.Cells.Item(3, 5).Value = 1.043512521010866
$ws.Cells.Item(3, 6).Value = 1.053918601910671
$ws.Cells.Item(3, 9).Value = 1.043245340974933
$ws.Cells.Item(3, 10).Value = 1.050660688316559
$ws.Cells.Item(3, 11).Value = 1.042726066261725
$ws.Cells.Item(3, 12).Value = 1.046099349593009
$ws.Cells.Item(3, 13).Value = 1.056478467404034
$ws.Cells.Item(3, 14).Value = 1.052152747006434

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.046950823501874
$ws.Cells.Item(4, 4).Value = 1.040601539992064
$ws.Cells.Item(4, 5).Value = 1.044379482102912
$ws.Cells.Item(4, 6).Value = 1.054995133406189
$ws.Cells.Item(4, 9).Value = 1.043537212512548
$ws.Cells.Item(4, 10).Value = 1.051429841727428
$ws.Cells.Item(4, 11).Value = 1.043073923294603
$ws.Cells.Item(4, 12).Value = 1.046842439025337
$ws.Cells.Item(4, 13).Value = 1.057431991774201
$ws.Cells.Item(4, 14).Value = 1.05292299270337

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.047370587420457
$ws.Cells.Item(5, 4).Value = 1.040799384556608
$ws.Cells.Item(5, 5).Value = 1.044743512984528
$ws.Cells.Item(5, 6).Value = 1.055447293671014
$ws.Cells.Item(5, 9).Value = 1.043659372220928
$ws.Cells.Item(5, 10).Value = 1.051752600288951
$ws.Cells.Item(5, 11).Value = 1.043219759099529
$ws.Cells.Item(5, 12).Value = 1.04715426685775
$ws.Cells.Item(5, 13).Value = 1.057832327656569
$ws.Cells.Item(5, 14).Value = 1.053246209619071

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.047441038301018
$ws.Cells.Item(6, 4).Value = 1.040832589273547
$ws.Cells.Item(6, 5).Value = 1.044804609823186
$ws.Cells.Item(6, 6).Value = 1.055523189424479
$ws.Cells.Item(6, 9).Value = 1.043679851615422
$ws.Cells.Item(6, 10).Value = 1.05180675828923
$ws.Cells.Item(6, 11).Value = 1.04324422199772
$ws.Cells.Item(6, 12).Value = 1.047206591086582
$ws.Cells.Item(6, 13).Value = 1.057899515261432
$ws.Cells.Item(6, 14).Value = 1.053300444529923

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.046956434380188
$ws.Cells.Item(7, 4).Value = 1.040604184557729
$ws.Cells.Item(7, 5).Value = 1.044384348018836
$ws.Cells.Item(7, 6).Value = 1.05500117680573
$ws.Cells.Item(7, 9).Value = 1.043538846946738
$ws.Cells.Item(7, 10).Value = 1.05143415676583
$ws.Cells.Item(7, 11).Value = 1.043075873541284
$ws.Cells.Item(7, 12).Value = 1.046846607902913
$ws.Cells.Item(7, 13).Value = 1.057437343135833
$ws.Cells.Item(7, 14).Value = 1.05292731386962

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.04492655571022
$ws.Cells.Item(8, 4).Value = 1.039647410956091
$ws.Cells.Item(8, 5).Value = 1.042623970794599
$ws.Cells.Item(8, 6).Value = 1.052815710958253
$ws.Cells.Item(8, 9).Value = 1.042944853201469
$ws.Cells.Item(8, 10).Value = 1.049871685498039
$ws.Cells.Item(8, 11).Value = 1.042368774991136
$ws.Cells.Item(8, 12).Value = 1.04533710580555
$ws.Cells.Item(8, 13).Value = 1.05550105125889
$ws.Cells.Item(8, 14).Value = 1.051362623713412

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.041332718762332
$ws.Cells.Item(9, 4).Value = 1.037953458494549
$ws.Cells.Item(9, 5).Value = 1.039507292146605
$ws.Cells.Item(9, 6).Value = 1.048950591216955
$ws.Cells.Item(9, 9).Value = 1.041880549308918
$ws.Cells.Item(9, 10).Value = 1.047098830583828
$ws.Cells.Item(9, 11).Value = 1.041109676980027
$ws.Cells.Item(9, 12).Value = 1.042658470975006
$ws.Cells.Item(9, 13).Value = 1.052071493295541
$ws.Cells.Item(9, 14).Value = 1.048585831027171

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.038924809534277
$ws.Cells.Item(10, 4).Value = 1.036818629733331
$ws.Cells.Item(10, 5).Value = 1.037419144955484
$ws.Cells.Item(10, 6).Value = 1.046363726707226
$ws.Cells.Item(10, 9).Value = 1.041158976983166
$ws.Cells.Item(10, 10).Value = 1.045236608095619
$ws.Cells.Item(10, 11).Value = 1.040261305987708
$ws.Cells.Item(10, 12).Value = 1.040859680743157
$ws.Cells.Item(10, 13).Value = 1.049772701022048
$ws.Cells.Item(10, 14).Value = 1.04672096396943

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.037879165876507
$ws.Cells.Item(11, 4).Value = 1.036325890701041
$ws.Cells.Item(11, 5).Value = 1.036512384800182
$ws.Cells.Item(11, 6).Value = 1.045241035714614
$ws.Cells.Item(11, 9).Value = 1.040843629206563
$ws.Cells.Item(11, 10).Value = 1.044426897745709
$ws.Cells.Item(11, 11).Value = 1.039891786545608
$ws.Cells.Item(11, 12).Value = 1.04007758989855
$ws.Cells.Item(11, 13).Value = 1.04877421591787
$ws.Cells.Item(11, 14).Value = 1.045910103737947

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.037490303776657
$ws.Cells.Item(12, 4).Value = 1.036142659472495
$ws.Cells.Item(12, 5).Value = 1.036175176573372
$ws.Cells.Item(12, 6).Value = 1.044823620417425
$ws.Cells.Item(12, 9).Value = 1.040726055102466
$ws.Cells.Item(12, 10).Value = 1.044125621602736
$ws.Cells.Item(12, 11).Value = 1.039754201376109
$ws.Cells.Item(12, 12).Value = 1.03978659635163
$ws.Cells.Item(12, 13).Value = 1.0484028575537
$ws.Cells.Item(12, 14).Value = 1.045608399748295

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.037573737189782
$ws.Cells.Item(13, 4).Value = 1.036181972565549
$ws.Cells.Item(13, 5).Value = 1.036247527010785
$ws.Cells.Item(13, 6).Value = 1.044913175630711
$ws.Cells.Item(13, 9).Value = 1.04075129511594
$ws.Cells.Item(13, 10).Value = 1.044190269756222
$ws.Cells.Item(13, 11).Value = 1.039783728812175
$ws.Cells.Item(13, 12).Value = 1.039849037782531
$ws.Cells.Item(13, 13).Value = 1.048482536939559
$ws.Cells.Item(13, 14).Value = 1.045673139709574

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.037847031940669
$ws.Cells.Item(14, 4).Value = 1.036310748971142
$ws.Cells.Item(14, 5).Value = 1.036484519220884
$ws.Cells.Item(14, 6).Value = 1.045206540228242
$ws.Cells.Item(14, 9).Value = 1.040833919490819
$ws.Cells.Item(14, 10).Value = 1.044402004699716
$ws.Cells.Item(14, 11).Value = 1.039880420447289
$ws.Cells.Item(14, 12).Value = 1.040053546333171
$ws.Cells.Item(14, 13).Value = 1.048743529109951
$ws.Cells.Item(14, 14).Value = 1.045885175340975

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.038015356069641
$ws.Cells.Item(15, 4).Value = 1.036390064972999
$ws.Cells.Item(15, 5).Value = 1.036630485027102
$ws.Cells.Item(15, 6).Value = 1.045387238725032
$ws.Cells.Item(15, 9).Value = 1.040884768687571
$ws.Cells.Item(15, 10).Value = 1.044532393234736
$ws.Cells.Item(15, 11).Value = 1.039939951633945
$ws.Cells.Item(15, 12).Value = 1.040179485575939
$ws.Cells.Item(15, 13).Value = 1.048904271496942
$ws.Cells.Item(15, 14).Value = 1.046015749042669

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.038994141140675
$ws.Cells.Item(16, 4).Value = 1.036851302480144
$ws.Cells.Item(16, 5).Value = 1.037479268484089
$ws.Cells.Item(16, 6).Value = 1.046438180932206
$ws.Cells.Item(16, 9).Value = 1.041179844100596
$ws.Cells.Item(16, 10).Value = 1.045290274394188
$ws.Cells.Item(16, 11).Value = 1.040285783792899
$ws.Cells.Item(16, 12).Value = 1.040911517291146
$ws.Cells.Item(16, 13).Value = 1.049838901016052
$ws.Cells.Item(16, 14).Value = 1.046774706480298

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.039607294493732
$ws.Cells.Item(17, 4).Value = 1.037140260798826
$ws.Cells.Item(17, 5).Value = 1.038010990713067
$ws.Cells.Item(17, 6).Value = 1.04709671424576
$ws.Cells.Item(17, 9).Value = 1.041364157386502
$ws.Cells.Item(17, 10).Value = 1.045764768116696
$ws.Cells.Item(17, 11).Value = 1.040502132096101
$ws.Cells.Item(17, 12).Value = 1.041369837637526
$ws.Cells.Item(17, 13).Value = 1.050424333024872
$ws.Cells.Item(17, 14).Value = 1.047249874038318

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.039964647555582
$ws.Cells.Item(18, 4).Value = 1.037308675186851
$ws.Cells.Item(18, 5).Value = 1.038320887153338
$ws.Cells.Item(18, 6).Value = 1.04748057931611
$ws.Cells.Item(18, 9).Value = 1.0414713844521
$ws.Cells.Item(18, 10).Value = 1.046041209145785
$ws.Cells.Item(18, 11).Value = 1.040628115329544
$ws.Cells.Item(18, 12).Value = 1.041636859729888
$ws.Cells.Item(18, 13).Value = 1.05076550800051
$ws.Cells.Item(18, 14).Value = 1.047526707645376

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.040086447080355
$ws.Cells.Item(19, 4).Value = 1.037366078164401
$ws.Cells.Item(19, 5).Value = 1.038426512034675
$ws.Cells.Item(19, 6).Value = 1.047611426005549
$ws.Cells.Item(19, 9).Value = 1.041507898775576
$ws.Cells.Item(19, 10).Value = 1.046135413870196
$ws.Cells.Item(19, 11).Value = 1.040671036993276
$ws.Cells.Item(19, 12).Value = 1.041727855321774
$ws.Cells.Item(19, 13).Value = 1.050881789727797
$ws.Cells.Item(19, 14).Value = 1.047621046151301

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.039541538914663
$ws.Cells.Item(20, 4).Value = 1.037109271786219
$ws.Cells.Item(20, 5).Value = 1.03795396768666
$ws.Cells.Item(20, 6).Value = 1.047026085420557
$ws.Cells.Item(20, 9).Value = 1.041344411279748
$ws.Cells.Item(20, 10).Value = 1.045713892922603
$ws.Cells.Item(20, 11).Value = 1.040478941619793
$ws.Cells.Item(20, 12).Value = 1.041320696167416
$ws.Cells.Item(20, 13).Value = 1.050361552579866
$ws.Cells.Item(20, 14).Value = 1.047198926595613

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.037766566342198
$ws.Cells.Item(21, 4).Value = 1.036272833241331
$ws.Cells.Item(21, 5).Value = 1.036414741963179
$ws.Cells.Item(21, 6).Value = 1.045120162754686
$ws.Cells.Item(21, 9).Value = 1.040809600846476
$ws.Cells.Item(21, 10).Value = 1.04433966829273
$ws.Cells.Item(21, 11).Value = 1.039851956284817
$ws.Cells.Item(21, 12).Value = 1.039993337266441
$ws.Cells.Item(21, 13).Value = 1.048666686671508
$ws.Cells.Item(21, 14).Value = 1.04582275040914

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.036647882857313
$ws.Cells.Item(22, 4).Value = 1.035745737642723
$ws.Cells.Item(22, 5).Value = 1.035444667096489
$ws.Cells.Item(22, 6).Value = 1.043919525008214
$ws.Cells.Item(22, 9).Value = 1.040470797230636
$ws.Cells.Item(22, 10).Value = 1.043472662032575
$ws.Cells.Item(22, 11).Value = 1.039455839432165
$ws.Cells.Item(22, 12).Value = 1.039155933465833
$ws.Cells.Item(22, 13).Value = 1.047598294138335
$ws.Cells.Item(22, 14).Value = 1.044954512900648

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.037241177048564
$ws.Cells.Item(23, 4).Value = 1.03602527516293
$ws.Cells.Item(23, 5).Value = 1.035959143662548
$ws.Cells.Item(23, 6).Value = 1.044556228995119
$ws.Cells.Item(23, 9).Value = 1.040650646175146
$ws.Cells.Item(23, 10).Value = 1.043932563808774
$ws.Cells.Item(23, 11).Value = 1.039666010275351
$ws.Cells.Item(23, 12).Value = 1.039600129344861
$ws.Cells.Item(23, 13).Value = 1.048164935132749
$ws.Cells.Item(23, 14).Value = 1.045415067790121

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.039571251923205
$ws.Cells.Item(24, 4).Value = 1.037123274789373
$ws.Cells.Item(24, 5).Value = 1.037979734704616
$ws.Cells.Item(24, 6).Value = 1.047058000307485
$ws.Cells.Item(24, 9).Value = 1.041353334559922
$ws.Cells.Item(24, 10).Value = 1.045736882231444
$ws.Cells.Item(24, 11).Value = 1.040489421043965
$ws.Cells.Item(24, 12).Value = 1.041342902035543
$ws.Cells.Item(24, 13).Value = 1.050389921280251
$ws.Cells.Item(24, 14).Value = 1.04722194855191

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.042263883167769
$ws.Cells.Item(25, 4).Value = 1.038392348746496
$ws.Cells.Item(25, 5).Value = 1.040314817718145
$ws.Cells.Item(25, 6).Value = 1.049951551080654
$ws.Cells.Item(25, 9).Value = 1.042157804398556
$ws.Cells.Item(25, 10).Value = 1.047818050391391
$ws.Cells.Item(25, 11).Value = 1.041436753239799
$ws.Cells.Item(25, 12).Value = 1.043353225485922
$ws.Cells.Item(25, 13).Value = 1.052960260271083
$ws.Cells.Item(25, 14).Value = 1.049306072209357
